$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'" + '70.423.00'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +5.25%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '3.612.49'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  +5.03%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'" + '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  +0.11%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '589.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  +3.39%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '190.51'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +3.07%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'" + '0.644'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  +1.66%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '3.602.14'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  +4.95%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.13%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '0.177'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  -0.36%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.660'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  +2.55%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '58.19'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  +5.09%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '0.0000291'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  +3.53%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '9.81'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +4.54%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '4.189.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = "'" + '3.612.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  +4.98%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = "'" + '19.37'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +4.42%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '70.315.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  +5.32%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'" + '12.49'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  +4.04%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '0.120'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +0.26%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  +4.06%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '491.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  +3.24%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '17.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  +16.93%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '5.37'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +8.02%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + '4.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  +6.32%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + '90.86'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  +1.42%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  +5.17%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  +1.42%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + '9.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  +6.33%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'" + '32.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  +2.75%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + '7.56'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  +8.84%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'" + '628.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  +6.57%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '12.26'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +5.36%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  +7.08%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'" + '65.35'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  +3.53%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '0.0₃0820'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  +6.18%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '38.11'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  +4.09%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '0.404'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  +3.38%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  +0.10%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  -1.24%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '3.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -0.50%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '3.304.74'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +5.61%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '3.09'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  +5.94%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '0.0445'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  +4.90%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '2.68'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  +1.95%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '3.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +2.89%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '0.137'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +1.94%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '9.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +4.39%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  -2.48%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +5.36%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  +0.14%  '
$ws.Range('E51').Style = 'Normal'
